$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text does not look like a pure number - safe to assign directly
$simpleUpdates = @{
    'D2' = '65.124.30'
    'E2' = '  -0.01%  '
    'D3' = '3.541.70'
    'E3' = '  +4.10%  '
    'E4' = '  -0.04%  '
    'E5' = '  +2.98%  '
    'E6' = '  +1.57%  '
    'D7' = '3.541.23'
    'E7' = '  +4.07%  '
    'E8' = '  +0.14%  '
    'E9' = '  +0.29%  '
    'E10' = '  +4.48%  '
    'E11' = '  -2.62%  '
    'E12' = '  +4.29%  '
    'D13' = '4.142.22'
    'E13' = '  +4.11%  '
    'E14' = '  +4.11%  '
    'E15' = '  +6.22%  '
    'D16' = '3.540.16'
    'E16' = '  +3.81%  '
    'E17' = '  +1.48%  '
    'D18' = '65.136.89'
    'E18' = '  +0.03%  '
    'E19' = '  +4.90%  '
    'E20' = '  +1.23%  '
    'E21' = '  +5.87%  '
    'E22' = '  +3.65%  '
    'E23' = '  +4.62%  '
    'D24' = '3.681.96'
    'E24' = '  +4.13%  '
    'E25' = '  +2.45%  '
    'E26' = '  +0.02%  '
    'E27' = '  +9.93%  '
    'E28' = '  +12.00%  '
    'E29' = '  +0.17%  '
    'E30' = '  +4.26%  '
    'E31' = '  +4.20%  '
    'D32' = '3.560.29'
    'E32' = '  +4.32%  '
    'E33' = '  +0.04%  '
    'E34' = '  +5.22%  '
    'B35' = 'Fetch.AI'
    'C35' = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    'E35' = '  +20.10%  '
    'B36' = 'Kaspa'
    'C36' = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    'E36' = '  +2.14%  '
    'E37' = '  +9.54%  '
    'E38' = '  +0.60%  '
    'E39' = '  +4.52%  '
    'E40' = '  +7.96%  '
    'E41' = '  +7.09%  '
    'E42' = '  +2.66%  '
    'E43' = '  +19.56%  '
    'E44' = '  -1.57%  '
    'E45' = '  -0.05%  '
    'E46' = '  +2.71%  '
    'E47' = '  +10.14%  '
    'E48' = '  +5.20%  '
    'E49' = '  +6.34%  '
    'D50' = '2.393.12'
    'E50' = '  +10.64%  '
    'E51' = '  +16.28%  '
}

# Cells whose new text looks like a number (e.g. '597.43', '1.00') -
# these must be forced to remain text so Excel does not convert them to numeric values
$textUpdates = @{
    'D5' = '597.43'
    'D6' = '138.05'
    'D9' = '0.494'
    'D11' = '6.92'
    'D12' = '0.386'
    'D14' = '0.0000183'
    'D15' = '27.34'
    'D19' = '10.03'
    'D20' = '5.89'
    'D21' = '14.25'
    'D22' = '392.72'
    'D23' = '0.572'
    'D26' = '0.999'
    'D28' = '7.81'
    'D29' = '1.00'
    'D31' = '8.29'
    'D34' = '23.85'
    'D35' = '1.36'
    'D36' = '0.144'
    'D37' = '1.58'
    'D38' = '170.01'
    'D40' = '5.00'
    'D41' = '0.0803'
    'D42' = '0.824'
    'D43' = '26.38'
    'D44' = '42.54'
    'D46' = '4.43'
    'D47' = '1.20'
    'D48' = '1.67'
    'D49' = '6.83'
    'D51' = '307.04'
}

foreach ($ref in $simpleUpdates.Keys) {
    $ws.Range($ref).Value = $simpleUpdates[$ref]
}

foreach ($ref in $textUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textUpdates[$ref]
    $cell.Style = "Normal"
}
